$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Rating/Year columns to remain text-typed (matches source t="inlineStr" cells)
# instead of being auto-converted to numbers by COM Value assignment.
$ws.Range("C2:D21").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = 'A Real Pain'
$ws.Cells.Item(2, 2).Value = 'Comedy'
$ws.Cells.Item(2, 3).Value = '7.4'
$ws.Cells.Item(2, 4).Value = '2024'
$ws.Cells.Item(3, 1).Value = 'Wicked'
$ws.Cells.Item(3, 2).Value = 'Fantasy'
$ws.Cells.Item(3, 3).Value = '8'
$ws.Cells.Item(3, 4).Value = '2024'
$ws.Cells.Item(4, 1).Value = 'Armor'
$ws.Cells.Item(4, 2).Value = 'Action'
$ws.Cells.Item(4, 3).Value = '3.5'
$ws.Cells.Item(4, 4).Value = '2024'
$ws.Cells.Item(5, 1).Value = 'A Different Man'
$ws.Cells.Item(5, 2).Value = 'Comedy'
$ws.Cells.Item(5, 3).Value = '7'
$ws.Cells.Item(5, 4).Value = '2024'
$ws.Cells.Item(6, 1).Value = '[HI] Yodha'
$ws.Cells.Item(6, 2).Value = 'Action'
$ws.Cells.Item(6, 3).Value = '5.7'
$ws.Cells.Item(6, 4).Value = '2024'
$ws.Cells.Item(7, 1).Value = '[HI] Singham Again'
$ws.Cells.Item(7, 2).Value = 'Action'
$ws.Cells.Item(7, 3).Value = '5.3'
$ws.Cells.Item(7, 4).Value = '2024'
$ws.Cells.Item(8, 1).Value = 'Saturday Night'
$ws.Cells.Item(8, 2).Value = 'Biography'
$ws.Cells.Item(8, 3).Value = '7'
$ws.Cells.Item(8, 4).Value = '2024'
$ws.Cells.Item(9, 1).Value = 'Mr. Monk''s Last Case: A Monk Movie'
$ws.Cells.Item(9, 2).Value = 'Action'
$ws.Cells.Item(9, 3).Value = '6.7'
$ws.Cells.Item(9, 4).Value = '2023'
$ws.Cells.Item(10, 1).Value = 'The Outrun'
$ws.Cells.Item(10, 2).Value = 'Action'
$ws.Cells.Item(10, 3).Value = '6.9'
$ws.Cells.Item(10, 4).Value = '2024'
$ws.Cells.Item(11, 1).Value = 'Wallace & Gromit: Vengeance Most Fowl'
$ws.Cells.Item(11, 2).Value = 'Adventure'
$ws.Cells.Item(11, 3).Value = '8'
$ws.Cells.Item(11, 4).Value = '2024'
$ws.Cells.Item(12, 1).Value = 'The Order'
$ws.Cells.Item(12, 2).Value = 'Crime'
$ws.Cells.Item(12, 3).Value = '7'
$ws.Cells.Item(12, 4).Value = '2024'
$ws.Cells.Item(13, 1).Value = 'Bird'
$ws.Cells.Item(13, 2).Value = 'Drama'
$ws.Cells.Item(13, 3).Value = '7.2'
$ws.Cells.Item(13, 4).Value = '2024'
$ws.Cells.Item(14, 1).Value = 'The Lord of the Rings: The War of the Rohirrim'
$ws.Cells.Item(14, 2).Value = 'Action'
$ws.Cells.Item(14, 3).Value = '6.6'
$ws.Cells.Item(14, 4).Value = '2024'
$ws.Cells.Item(15, 1).Value = '[DA] The Promised Land'
$ws.Cells.Item(15, 2).Value = 'Action'
$ws.Cells.Item(15, 3).Value = '7.7'
$ws.Cells.Item(15, 4).Value = '2023'
$ws.Cells.Item(16, 1).Value = 'The Six Triple Eight'
$ws.Cells.Item(16, 2).Value = 'Drama'
$ws.Cells.Item(16, 3).Value = '6.5'
$ws.Cells.Item(16, 4).Value = '2024'
$ws.Cells.Item(17, 1).Value = 'Gladiator II'
$ws.Cells.Item(17, 2).Value = 'Action'
$ws.Cells.Item(17, 3).Value = '6.8'
$ws.Cells.Item(17, 4).Value = '2024'
$ws.Cells.Item(18, 1).Value = 'DragonHeart'
$ws.Cells.Item(18, 2).Value = 'Action'
$ws.Cells.Item(18, 3).Value = '6.4'
$ws.Cells.Item(18, 4).Value = '1996'
$ws.Cells.Item(19, 1).Value = 'Megalopolis'
$ws.Cells.Item(19, 2).Value = 'Drama'
$ws.Cells.Item(19, 3).Value = '4.8'
$ws.Cells.Item(19, 4).Value = '2024'
$ws.Cells.Item(20, 1).Value = 'Hush'
$ws.Cells.Item(20, 2).Value = 'Action'
$ws.Cells.Item(20, 3).Value = '6.6'
$ws.Cells.Item(20, 4).Value = '2016'
$ws.Cells.Item(21, 1).Value = 'Candyman'
$ws.Cells.Item(21, 2).Value = 'Action'
$ws.Cells.Item(21, 3).Value = '6.7'
$ws.Cells.Item(21, 4).Value = '1992'

# Restore the default cell style so only the values differ (NumberFormat above
# was only a scratch trick to stop auto-numeric coercion).
$ws.Range("C2:D21").Style = "Normal"
